$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 85.83574689470727

$ws.Range("N2").Value = $newValue
$ws.Range("N3").Value = $newValue
$ws.Range("N4").Value = $newValue
$ws.Range("N5").Value = $newValue
